# -----------------------------------------------------------------------
# Applies the "tambah coding di bab 2 dan tambah json di dftr isi" edit:
#   1. Removes the stray _GoBack bookmark that used to sit in the
#      "DAFTAR DAFTAR GAMBAR ... viii" TOC line.
#   2. Inserts a new "JSON ... 23" Daftar Isi entry (with the same
#      numbering / tab formatting as the existing "Android" entry) right
#      before the "Android" entry, and leaves the _GoBack bookmark at the
#      end of that freshly typed line (mirroring where Word drops it
#      after the last edit).
#   3. Moves the w:lastRenderedPageBreak hint from the "Halaman" running
#      header onto the "Jadwal Kegiatan" TOC heading above it (page break
#      now falls a line earlier because of the extra TOC entry above).
#   4. Adds the same w:lastRenderedPageBreak hint to "LAMPIRAN LAMPIRAN".
# -----------------------------------------------------------------------

$d = $word.ActiveDocument

$wNs = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

# ---- helpers ----------------------------------------------------------

function Find-ParagraphIndex($doc, $text) {
    $i = 1
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Contains($text)) {
            return $i
        }
        $i = $i + 1
    }
    return -1
}

function Get-ParaXml($para) {
    $full = $para.Range.WordOpenXML
    $m = [regex]::Match($full, '(?s)<w:p[ >].*?</w:p>')
    return $m.Value
}

# =========================================================================
# 1) Drop the old _GoBack bookmark (was sitting inside the "DAFTAR DAFTAR
#    GAMBAR" TOC line, between the 3rd and 4th roman-numeral run).
# =========================================================================
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# =========================================================================
# 2) Insert the new "JSON ... 23" entry right before "Android" in the
#    Daftar Isi, re-using the exact paragraph formatting of the existing
#    "Android" line (pStyle ListParagraph, ilvl 2, numId 4, dotted tab
#    leader, ind left 709).
# =========================================================================
$androidIdx = Find-ParagraphIndex $d "Android"
$androidPara = $d.Paragraphs($androidIdx)
$null = $androidPara.Range.InsertParagraphBefore()

$jsonPara = $d.Paragraphs($androidIdx)
$jsonXml = "<w:p xmlns:w='$wNs'>" +
    "<w:pPr>" +
        "<w:pStyle w:val='ListParagraph'/>" +
        "<w:numPr><w:ilvl w:val='2'/><w:numId w:val='4'/></w:numPr>" +
        "<w:tabs><w:tab w:val='left' w:leader='dot' w:pos='7088'/><w:tab w:val='right' w:pos='7938'/></w:tabs>" +
        "<w:spacing w:after='0'/>" +
        "<w:ind w:left='709'/>" +
        "<w:rPr><w:rFonts w:eastAsia='Times New Roman' w:cs='Times New Roman'/><w:szCs w:val='24'/><w:lang w:eastAsia='id-ID'/></w:rPr>" +
    "</w:pPr>" +
    "<w:r><w:rPr><w:rFonts w:eastAsia='Times New Roman' w:cs='Times New Roman'/><w:szCs w:val='24'/><w:lang w:val='en-ID' w:eastAsia='id-ID'/></w:rPr><w:t>JSON</w:t></w:r>" +
    "<w:r><w:rPr><w:rFonts w:eastAsia='Times New Roman' w:cs='Times New Roman'/><w:szCs w:val='24'/><w:lang w:val='en-ID' w:eastAsia='id-ID'/></w:rPr><w:tab/></w:r>" +
    "<w:r><w:rPr><w:rFonts w:eastAsia='Times New Roman' w:cs='Times New Roman'/><w:szCs w:val='24'/><w:lang w:val='en-ID' w:eastAsia='id-ID'/></w:rPr><w:tab/><w:t>23</w:t></w:r>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
    "</w:p>"
$null = $jsonPara.Range.InsertXML($jsonXml)

# =========================================================================
# 3) Move w:lastRenderedPageBreak from "Halaman" onto "Jadwal Kegiatan".
# =========================================================================
$jadwalIdx = Find-ParagraphIndex $d "Jadwal Kegiatan"
$jadwalPara = $d.Paragraphs($jadwalIdx)
$jadwalXml = Get-ParaXml $jadwalPara
$jadwalXml = $jadwalXml.Replace('<w:t>Jadwal Kegiatan</w:t>', '<w:lastRenderedPageBreak/><w:t>Jadwal Kegiatan</w:t>')
$null = $jadwalPara.Range.InsertXML($jadwalXml)

$halamanIdx = $jadwalIdx + 1
$halamanPara = $d.Paragraphs($halamanIdx)
$halamanXml = Get-ParaXml $halamanPara
$halamanXml = $halamanXml.Replace('<w:lastRenderedPageBreak/>', '')
$null = $halamanPara.Range.InsertXML($halamanXml)

# =========================================================================
# 4) Add w:lastRenderedPageBreak to "LAMPIRAN LAMPIRAN".
# =========================================================================
$lampiranIdx = Find-ParagraphIndex $d "LAMPIRAN LAMPIRAN"
$lampiranPara = $d.Paragraphs($lampiranIdx)
$lampiranXml = Get-ParaXml $lampiranPara
$lampiranXml = $lampiranXml.Replace('<w:t>LAMPIRAN LAMPIRAN</w:t>', '<w:lastRenderedPageBreak/><w:t>LAMPIRAN LAMPIRAN</w:t>')
# WordOpenXML round-tripping drops the run-level rsid on the "58" run;
# put it back so this paragraph comes out byte-identical apart from the
# intended w:lastRenderedPageBreak addition.
$lampiranXml = $lampiranXml.Replace('<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:bCs/><w:szCs w:val="24"/><w:lang w:val="en-ID" w:eastAsia="id-ID"/></w:rPr><w:t>58</w:t></w:r>', '<w:r w:rsidR="003A0AF3"><w:rPr><w:rFonts w:eastAsia="Times New Roman" w:cs="Times New Roman"/><w:bCs/><w:szCs w:val="24"/><w:lang w:val="en-ID" w:eastAsia="id-ID"/></w:rPr><w:t>58</w:t></w:r>')
$null = $lampiranPara.Range.InsertXML($lampiranXml)

Write-Output "done"
